$d = $word.ActiveDocument

# Commit: "Se quito 'capacidad' de Habitacion y se modifico eso en todo lado"
# The "Habitacion" bullet in the entity-attribute list loses its "capacidad"
# attribute, so "Habitación: id, capacidad, disponibilidad" becomes
# "Habitación: id, disponible".
$d.Content.Find.Execute(
    "Habitación: id, capacidad, disponibilidad",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Habitación: id, disponible",
    2
) | Out-Null
